# Refresh the latest crypto price/volume snapshot (rows 2-51, columns D "Price"
# and E "Volume(1h)") to the values captured by the Tue May 14 14:40:03 UTC 2024
# GitHub Actions run. Columns A (id), B (Coin) and C (Link) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '61.936.55'; E = '  -1.55%  '; Quote = $false },
    @{ Row = 3; D = '2.909.45'; E = '  -1.78%  '; Quote = $false },
    @{ Row = 4; D = '1.00'; E = '  +0.03%  '; Quote = $true },
    @{ Row = 5; D = '572.86'; E = '  -3.70%  '; Quote = $true },
    @{ Row = 6; D = '144.52'; E = '  -1.22%  '; Quote = $true },
    @{ Row = 7; D = $null; E = '  +0.20%  '; Quote = $false },
    @{ Row = 8; D = $null; E = '  -0.82%  '; Quote = $false },
    @{ Row = 9; D = '2.904.75'; E = '  -1.96%  '; Quote = $false },
    @{ Row = 10; D = '6.68'; E = '  -7.72%  '; Quote = $true },
    @{ Row = 11; D = $null; E = '  -1.16%  '; Quote = $false },
    @{ Row = 12; D = '0.433'; E = '  -2.69%  '; Quote = $true },
    @{ Row = 13; D = '0.0000233'; E = '  -2.90%  '; Quote = $true },
    @{ Row = 14; D = '32.26'; E = '  -2.98%  '; Quote = $true },
    @{ Row = 15; D = $null; E = '  -0.75%  '; Quote = $false },
    @{ Row = 16; D = '3.393.70'; E = '  -1.78%  '; Quote = $false },
    @{ Row = 17; D = '61.961.38'; E = '  -1.40%  '; Quote = $false },
    @{ Row = 18; D = '6.63'; E = $null; Quote = $true },
    @{ Row = 19; D = '2.913.31'; E = '  -1.22%  '; Quote = $false },
    @{ Row = 20; D = '436.36'; E = '  -1.33%  '; Quote = $true },
    @{ Row = 21; D = '13.27'; E = '  -2.14%  '; Quote = $true },
    @{ Row = 22; D = '0.657'; E = '  -2.24%  '; Quote = $true },
    @{ Row = 23; D = '6.91'; E = '  -2.74%  '; Quote = $true },
    @{ Row = 24; D = '79.33'; E = '  -2.79%  '; Quote = $true },
    @{ Row = 25; D = '11.99'; E = '  +1.03%  '; Quote = $true },
    @{ Row = 26; D = $null; E = '  -9.53%  '; Quote = $false },
    @{ Row = 27; D = $null; E = '  +0.04%  '; Quote = $false },
    @{ Row = 28; D = $null; E = '  -4.50%  '; Quote = $false },
    @{ Row = 29; D = '0.0000109'; E = '  +12.48%  '; Quote = $true },
    @{ Row = 30; D = '7.08'; E = '  -2.71%  '; Quote = $true },
    @{ Row = 31; D = '2.53'; E = '  -3.33%  '; Quote = $true },
    @{ Row = 32; D = '2.07'; E = '  -3.90%  '; Quote = $true },
    @{ Row = 33; D = $null; E = '  -2.39%  '; Quote = $false },
    @{ Row = 34; D = '1.00'; E = '  -0.10%  '; Quote = $true },
    @{ Row = 35; D = '25.67'; E = '  -3.32%  '; Quote = $true },
    @{ Row = 36; D = '0.959'; E = '  -3.70%  '; Quote = $true },
    @{ Row = 37; D = $null; E = '  -3.65%  '; Quote = $false },
    @{ Row = 38; D = '2.96'; E = '  -4.76%  '; Quote = $true },
    @{ Row = 39; D = '49.14'; E = '  -0.88%  '; Quote = $true },
    @{ Row = 40; D = $null; E = '  -3.90%  '; Quote = $false },
    @{ Row = 41; D = $null; E = '  -1.12%  '; Quote = $false },
    @{ Row = 42; D = '8.27'; E = '  -3.03%  '; Quote = $true },
    @{ Row = 43; D = '0.270'; E = '  -4.30%  '; Quote = $true },
    @{ Row = 44; D = '38.90'; E = '  -4.93%  '; Quote = $true },
    @{ Row = 45; D = '2.688.43'; E = '  -1.97%  '; Quote = $false },
    @{ Row = 46; D = '133.40'; E = '  -0.63%  '; Quote = $true },
    @{ Row = 47; D = $null; E = '  -1.18%  '; Quote = $false },
    @{ Row = 48; D = $null; E = '  +0.03%  '; Quote = $false },
    @{ Row = 49; D = '340.18'; E = '  -7.14%  '; Quote = $true },
    @{ Row = 50; D = $null; E = '  -2.03%  '; Quote = $false },
    @{ Row = 51; D = '21.72'; E = '  -5.63%  '; Quote = $true }
)

foreach ($u in $updates) {
    $dCell = $ws.Range("D" + $u.Row)
    $eCell = $ws.Range("E" + $u.Row)

    if ($null -ne $u.D) {
        if ($u.Quote) {
            # Values such as "1.00" or "6.68" would otherwise be auto-coerced to
            # numbers by Excel's type inference; force text storage (General
            # formatting, matching the original cells) without leaving a visible
            # quote-prefix style behind.
            $dCell.NumberFormat = "@"
            $dCell.Value = $u.D
            $dCell.Style = "Normal"
        } else {
            $dCell.Value = $u.D
        }
    }

    if ($null -ne $u.E) {
        $eCell.Value = $u.E
    }
}
